$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(334, 2).Value = 0.8810912419855594
$ws.Cells.Item(335, 2).Value = 0.8214963534474373
$ws.Cells.Item(336, 2).Value = 0.8869845113158226
$ws.Cells.Item(337, 2).Value = 0.8105607953667641
$ws.Cells.Item(338, 2).Value = 0.7763261064887047
$ws.Cells.Item(339, 2).Value = 1.029825215339661
$ws.Cells.Item(340, 2).Value = 0.7396623514592647
$ws.Cells.Item(341, 2).Value = 0.9300873246788979
$ws.Cells.Item(342, 2).Value = 0.7431439656019211
$ws.Cells.Item(343, 2).Value = 1.252684499621391
$ws.Cells.Item(344, 2).Value = 0.5400192460417748
$ws.Cells.Item(345, 2).Value = 0.4820150673389435
$ws.Cells.Item(346, 2).Value = 0.102651838362217
$ws.Cells.Item(347, 2).Value = 0.5836758437752724
$ws.Cells.Item(348, 2).Value = 0.3101505497097969
$ws.Cells.Item(349, 2).Value = 0.4494251096248627
$ws.Cells.Item(350, 2).Value = 0.3493212205171585
$ws.Cells.Item(351, 2).Value = -0.07923666030168532
$ws.Cells.Item(352, 2).Value = 1.01478587731719
$ws.Cells.Item(353, 2).Value = 0.3373434242606163
$ws.Cells.Item(354, 2).Value = 0.580757271796465
$ws.Cells.Item(355, 2).Value = 0.05624767482280732
$ws.Cells.Item(356, 2).Value = 0.3945392292737961
$ws.Cells.Item(357, 2).Value = 0.4908187028765679
$ws.Cells.Item(358, 2).Value = 0.4117062291502953
$ws.Cells.Item(359, 2).Value = 0.3670210695266723
$ws.Cells.Item(360, 2).Value = -0.2127494010329246
$ws.Cells.Item(361, 2).Value = 0.3056404262781143
$ws.Cells.Item(362, 2).Value = 0.9879850387573242
$ws.Cells.Item(363, 2).Value = 0.4814973147213459
$ws.Cells.Item(364, 2).Value = 0.3756522786617279
$ws.Cells.Item(365, 2).Value = 0.8622520268708468
$ws.Cells.Item(366, 2).Value = 0.4422421139478683
$ws.Cells.Item(367, 2).Value = 0.3255287966132164
$ws.Cells.Item(368, 2).Value = 0.2913089412450791
$ws.Cells.Item(369, 2).Value = 0.1962499183416367
